# Insert a new price-report row at row 28 (pushing existing rows 28..97
# down to 29..98) and populate it with the new "New Hall" orange entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(28).Insert()

$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44791
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100102
$ws.Range("H28").Value = "Cítricos"
$ws.Range("I28").Value = 100102005
$ws.Range("J28").Value = "Naranja"
$ws.Range("K28").Value = "New Hall"
$ws.Range("L28").Value = "Tercera"
$ws.Range("M28").Value = 300
$ws.Range("N28").Value = 500
$ws.Range("O28").Value = 600
$ws.Range("P28").Value = 550
$ws.Range("Q28").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 550
$ws.Range("T28").Value = 1
